$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that previously only had columns A-F populated (F = "N/A" shared
# string) now get "N/A" filled across columns G:N as well, matching the
# existing F-column value for that row.
$ws.Range("G2:N2").Value = "N/A"
$ws.Range("G4:N4").Value = "N/A"
$ws.Range("G6:N6").Value = "N/A"
$ws.Range("G12:N12").Value = "N/A"
$ws.Range("G19:N19").Value = "N/A"
$ws.Range("G20:N20").Value = "N/A"
$ws.Range("G22:N22").Value = "N/A"
$ws.Range("G26:N26").Value = "N/A"
$ws.Range("G29:N29").Value = "N/A"
$ws.Range("G34:N34").Value = "N/A"
$ws.Range("G37:N37").Value = "N/A"
$ws.Range("G39:N39").Value = "N/A"
$ws.Range("G42:N42").Value = "N/A"
$ws.Range("G44:N44").Value = "N/A"
$ws.Range("G48:N48").Value = "N/A"
$ws.Range("G49:N49").Value = "N/A"
$ws.Range("G52:N52").Value = "N/A"
$ws.Range("G53:N53").Value = "N/A"
$ws.Range("G56:N56").Value = "N/A"
$ws.Range("G59:N59").Value = "N/A"
$ws.Range("G62:N62").Value = "N/A"
$ws.Range("G66:N66").Value = "N/A"
$ws.Range("G67:N67").Value = "N/A"
$ws.Range("G70:N70").Value = "N/A"
$ws.Range("G74:N74").Value = "N/A"
$ws.Range("G75:N75").Value = "N/A"
$ws.Range("G79:N79").Value = "N/A"
$ws.Range("G81:N81").Value = "N/A"
$ws.Range("G82:N82").Value = "N/A"
$ws.Range("G85:N85").Value = "N/A"
$ws.Range("G87:N87").Value = "N/A"
$ws.Range("G94:N94").Value = "N/A"
$ws.Range("G95:N95").Value = "N/A"
$ws.Range("G97:N97").Value = "N/A"
$ws.Range("G98:N98").Value = "N/A"
$ws.Range("G101:N101").Value = "N/A"
$ws.Range("G104:N104").Value = "N/A"
$ws.Range("G107:N107").Value = "N/A"
$ws.Range("G108:N108").Value = "N/A"
$ws.Range("G109:N109").Value = "N/A"
$ws.Range("G111:N111").Value = "N/A"
$ws.Range("G113:N113").Value = "N/A"
$ws.Range("G115:N115").Value = "N/A"
$ws.Range("G118:N118").Value = "N/A"
$ws.Range("G121:N121").Value = "N/A"
$ws.Range("G122:N122").Value = "N/A"
$ws.Range("G126:N126").Value = "N/A"
$ws.Range("G127:N127").Value = "N/A"
$ws.Range("G133:N133").Value = "N/A"
$ws.Range("G134:N134").Value = "N/A"
$ws.Range("G138:N138").Value = "N/A"
$ws.Range("G142:N142").Value = "N/A"
$ws.Range("G144:N144").Value = "N/A"
$ws.Range("G151:N151").Value = "N/A"
$ws.Range("G156:N156").Value = "N/A"
$ws.Range("G159:N159").Value = "N/A"
$ws.Range("G165:N165").Value = "N/A"
$ws.Range("G168:N168").Value = "N/A"
$ws.Range("G170:N170").Value = "N/A"
$ws.Range("G174:N174").Value = "N/A"
$ws.Range("G179:N179").Value = "N/A"
$ws.Range("G184:N184").Value = "N/A"
$ws.Range("G187:N187").Value = "N/A"
$ws.Range("G189:N189").Value = "N/A"
